$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (G..K / cols 7..11) must be forced to Text format so that
# numeric-looking strings (e.g. '0', '168.75') are stored as text, not numbers.
$ws.Range("G1:K8").NumberFormat = "@"

# Header row
$ws.Cells.Item(1, 1).Value = "venue"
$ws.Cells.Item(1, 2).Value = "date"
$ws.Cells.Item(1, 3).Value = "result"
$ws.Cells.Item(1, 4).Value = "ownTeam"
$ws.Cells.Item(1, 5).Value = "oppTeam"
$ws.Cells.Item(1, 6).Value = "batsman"
$ws.Cells.Item(1, 7).Value = "totalRuns"
$ws.Cells.Item(1, 8).Value = "totalBalls"
$ws.Cells.Item(1, 9).Value = "total4s"
$ws.Cells.Item(1, 10).Value = "total6s"
$ws.Cells.Item(1, 11).Value = "sr"

# Row 2
$ws.Cells.Item(2, 1).Value = " Abu Dhabi"
$ws.Cells.Item(2, 2).Value = " October 30 2020"
$ws.Cells.Item(2, 3).Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Cells.Item(2, 4).Value = "Kings XI Punjab"
$ws.Cells.Item(2, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(2, 6).Value = "Mandeep Singh "
$ws.Cells.Item(2, 7).Value = "0"
$ws.Cells.Item(2, 8).Value = "1"
$ws.Cells.Item(2, 9).Value = "0"
$ws.Cells.Item(2, 10).Value = "0"
$ws.Cells.Item(2, 11).Value = "0.00"

# Row 3
$ws.Cells.Item(3, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(3, 2).Value = " October 04 2020"
$ws.Cells.Item(3, 3).Value = "Super Kings won by 10 wickets (with 14 balls remaining)"
$ws.Cells.Item(3, 4).Value = "Kings XI Punjab"
$ws.Cells.Item(3, 5).Value = "Chennai Super Kings"
$ws.Cells.Item(3, 6).Value = "Mandeep Singh "
$ws.Cells.Item(3, 7).Value = "27"
$ws.Cells.Item(3, 8).Value = "16"
$ws.Cells.Item(3, 9).Value = "0"
$ws.Cells.Item(3, 10).Value = "2"
$ws.Cells.Item(3, 11).Value = "168.75"

# Row 4
$ws.Cells.Item(4, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(4, 2).Value = " October 24 2020"
$ws.Cells.Item(4, 3).Value = "Kings XI won by 12 runs"
$ws.Cells.Item(4, 4).Value = "Kings XI Punjab"
$ws.Cells.Item(4, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(4, 6).Value = "Mandeep Singh "
$ws.Cells.Item(4, 7).Value = "17"
$ws.Cells.Item(4, 8).Value = "14"
$ws.Cells.Item(4, 9).Value = "1"
$ws.Cells.Item(4, 10).Value = "0"
$ws.Cells.Item(4, 11).Value = "121.42"

# Row 5
$ws.Cells.Item(5, 1).Value = " Abu Dhabi"
$ws.Cells.Item(5, 2).Value = " November 01 2020"
$ws.Cells.Item(5, 3).Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Cells.Item(5, 4).Value = "Kings XI Punjab"
$ws.Cells.Item(5, 5).Value = "Chennai Super Kings"
$ws.Cells.Item(5, 6).Value = "Mandeep Singh "
$ws.Cells.Item(5, 7).Value = "14"
$ws.Cells.Item(5, 8).Value = "15"
$ws.Cells.Item(5, 9).Value = "1"
$ws.Cells.Item(5, 10).Value = "0"
$ws.Cells.Item(5, 11).Value = "93.33"

# Row 6
$ws.Cells.Item(6, 1).Value = " Abu Dhabi"
$ws.Cells.Item(6, 2).Value = " October 10 2020"
$ws.Cells.Item(6, 3).Value = "KKR won by 2 runs"
$ws.Cells.Item(6, 4).Value = "Kings XI Punjab"
$ws.Cells.Item(6, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(6, 6).Value = "Mandeep Singh "
$ws.Cells.Item(6, 7).Value = "0"
$ws.Cells.Item(6, 8).Value = "1"
$ws.Cells.Item(6, 9).Value = "0"
$ws.Cells.Item(6, 10).Value = "0"
$ws.Cells.Item(6, 11).Value = "0.00"

# Row 7
$ws.Cells.Item(7, 1).Value = " Sharjah"
$ws.Cells.Item(7, 2).Value = " October 26 2020"
$ws.Cells.Item(7, 3).Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Cells.Item(7, 4).Value = "Kings XI Punjab"
$ws.Cells.Item(7, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(7, 6).Value = "Mandeep Singh "
$ws.Cells.Item(7, 7).Value = "66"
$ws.Cells.Item(7, 8).Value = "56"
$ws.Cells.Item(7, 9).Value = "8"
$ws.Cells.Item(7, 10).Value = "2"
$ws.Cells.Item(7, 11).Value = "117.85"

# Row 8
$ws.Cells.Item(8, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(8, 2).Value = " October 08 2020"
$ws.Cells.Item(8, 3).Value = "Sunrisers won by 69 runs"
$ws.Cells.Item(8, 4).Value = "Kings XI Punjab"
$ws.Cells.Item(8, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(8, 6).Value = "Mandeep Singh "
$ws.Cells.Item(8, 7).Value = "6"
$ws.Cells.Item(8, 8).Value = "6"
$ws.Cells.Item(8, 9).Value = "0"
$ws.Cells.Item(8, 10).Value = "0"
$ws.Cells.Item(8, 11).Value = "100.00"
